$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(46009, 10528.2037874986, 9506.96617785098, 11232.26, 6633.91412668432, 204.525846022304)
    3  = @(46010, 10578.8430243892, 8945.32849587728, 11232.26, 6920.85273091918, 193.080051116519)
    4  = @(46011, 4200.00459288774, 6358.5000290554, 11232.26, 6802.67171624156, 80.3713227207064)
    5  = @(46012, 4228.95362472504, 6499.10652903338, 11232.26, 7003.28932063639, 94.5889937362404)
    6  = @(46013, 10055.3141109786, 9474.25373934299, 11232.26, 7480.43738797911, 238.434630305088)
    7  = @(46014, 9566.18303866457, 9369.61042170503, 17376.26, 8113.18340118371, 4.43890928703104)
    8  = @(46015, 9566.18303866457, 9550.27057415864, 17208.26, 8113.18340118371, 18.9664156392647)
    9  = @(46016, 8471.80799626025, 8913.38728242721, 16416.26, 7768.80692669759, 11.0805920468664)
    10 = @(46017, 8471.80799626025, 8143.56085254162, 15840.26, 7768.80692669759, 3.00449080163366)
    11 = @(46018, 8547.24893097507, 8729.58716086958, 16632.26, 7775.98586474824, -5.27862393259056)
    12 = @(46019, 8471.80799626025, 8896.14257139484, 17448.26, 7768.80692669759, -32.6379375794821)
    13 = @(46020, 9603.52690270833, 10485.5893446834, 19176.26, 8403.10977092558, -11.981703516291)
    14 = @(46021, 9603.52690270833, 10905.6554140192, 11232.26, 8403.10977092558, 336.521049372701)
    15 = @(46022, 9603.52690270833, 10819.0659278073, 11232.26, 8403.10977092558, 332.91315411387)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
}
